# Automatische test-sync: 2025-06-19 17:59:30
#
# Appends the newly received e-mail ("Offerte voor zakelijke samenwerking")
# to the "Logs" sheet as row 33, extends the conditional formatting ranges
# to cover the new row, and updates the "Dashboard" summary sheet so the
# "Offerte-aanvraag" category count reflects the new entry (1 -> 2), while
# keeping the category list ordered by descending count.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Logs sheet: append the new mail as row 33
# ---------------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A33").Value = "Offerte voor zakelijke samenwerking"
$logs.Range("B33").Value = "mailmind.test@zohomail.eu"
$logs.Range("C33").Value = "Kunt u mij een offerte sturen voor 100 stuks product X?"
$logs.Range("D33").Value = "Offerte-aanvraag"
$logs.Range("F33").Value = "2025-06-19 17:59:20"
$logs.Range("G33").Value = "Nee"

# ---------------------------------------------------------------------------
# 2) Logs sheet: extend conditional formatting ranges to include row 33
# ---------------------------------------------------------------------------
$categoryRules = $logs.Range("D2:D32").FormatConditions
$categoryRules.Item(1).ModifyAppliesToRange($logs.Range("D2:D33"))

$answeredRules = $logs.Range("G2:G32").FormatConditions
$answeredRules.Item(1).ModifyAppliesToRange($logs.Range("G2:G33"))

# ---------------------------------------------------------------------------
# 3) Dashboard sheet: "Offerte-aanvraag" now has 2 occurrences, which moves
#    it above "Informatieaanvraag" (still 1) in the ranked category table.
# ---------------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A7").Value = "Offerte-aanvraag"
$dashboard.Range("B7").Value = 2
$dashboard.Range("A8").Value = "Informatieaanvraag"
$dashboard.Range("B8").Value = 1
